# Update countries & provincias Spain
# Applies the data refresh captured in the commit diff:
#  - swaps the country labels for a few tied/re-ordered rows
#  - refreshes the numeric COVID figures for Honduras, Haiti/Tayikistan,
#    Mongolia and Islas Turcas y Caicos / San Martin (Parte Francesa)
#  - bumps the "last updated" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 3 de Julio de 2020 a las 05:42"

# --- Row 55: Honduras, refreshed figures ---
$ws.Range("B55").Value = 21120
$ws.Range("C55").Value = 858
$ws.Range("D55").Value = 2190
$ws.Range("E55").Value = 18339
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 49
$ws.Range("H55").Value = 591

# --- Rows 84/85: Haiti and Tayikistan swap places with refreshed data ---
# Row 84 now shows Haiti with brand-new figures
$ws.Range("A84").Value = "Haiti"
$ws.Range("B84").Value = 6101
$ws.Range("C84").Value = 61
$ws.Range("D84").Value = 1141
$ws.Range("E84").Value = 4850
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 3
$ws.Range("H84").Value = 110

# Row 85 now shows Tayikistan, carrying the figures Haiti used to have
$ws.Range("A85").Value = "Tayikistan"
$ws.Range("B85").Value = 6058
$ws.Range("C85").Value = 0
$ws.Range("D85").Value = 4690
$ws.Range("E85").Value = 1316
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 52

# --- Row 169: Mongolia, refreshed figures ---
$ws.Range("D169").Value = 179
$ws.Range("E169").Value = 41

# --- Rows 203/204: Santa Lucia <-> Laos (tied values, label swap only) ---
$ws.Range("A203").Value = "Laos"
$ws.Range("A204").Value = "Santa Lucia"

# --- Rows 205/206: Fiyi <-> Dominica (tied values, label swap only) ---
$ws.Range("A205").Value = "Dominica"
$ws.Range("A206").Value = "Fiyi"

# --- Rows 192/193: Islas Turcas y Caicos and San Martin swap with refreshed data ---
# Row 192 now shows Islas Turcas y Caicos with brand-new figures
$ws.Range("A192").Value = "Islas Turcas y Caicos"
$ws.Range("B192").Value = 44
$ws.Range("C192").Value = 2
$ws.Range("D192").Value = 11
$ws.Range("E192").Value = 31
$ws.Range("F192").Value = 0
$ws.Range("G192").Value = 0
$ws.Range("H192").Value = 2

# Row 193 now shows San Martin (Parte Francesa), carrying the figures
# Islas Turcas y Caicos used to have
$ws.Range("A193").Value = "San Martin (Parte Francesa)"
$ws.Range("B193").Value = 43
$ws.Range("C193").Value = 0
$ws.Range("D193").Value = 37
$ws.Range("E193").Value = 3
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 3

# --- Rows 209/210: Groenlandia <-> Islas Malvinas (tied values, label swap only) ---
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"
